$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 334; everything from old row 334 onward shifts down by one.
$ws.Rows.Item(334).Insert()

# Populate the newly inserted row 334 with the new record's data.
$ws.Cells.Item(334, 1).Value = 10
$ws.Cells.Item(334, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(334, 3).Value = "La Araucanía"
$ws.Cells.Item(334, 4).Value = 44714
$ws.Cells.Item(334, 5).Value = 9
$ws.Cells.Item(334, 6).Value = 100112040
$ws.Cells.Item(334, 7).Value = "Cilantro"
$ws.Cells.Item(334, 8).Value = "Sin especificar"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 65
$ws.Cells.Item(334, 11).Value = 6000
$ws.Cells.Item(334, 12).Value = 6000
$ws.Cells.Item(334, 13).Value = 6000
$ws.Cells.Item(334, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(334, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(334, 16).Value = 3000
$ws.Cells.Item(334, 17).Value = 2
$ws.Cells.Item(334, 18).Value = "Hortaliza"
